$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B6:E6").Delete()
